$ws = $excel.ActiveWorkbook.ActiveSheet

$kValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 0
    6 = 2
    7 = 0
    8 = 0
    9 = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 2
    19 = 0
    20 = 3
    21 = 1
    22 = 3
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 2
    44 = 1
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 3
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 2
    55 = 2
    56 = 0
    57 = 2
    58 = 1
    59 = 2
    60 = 0
    61 = 1
    62 = 2
    63 = 0
    64 = 2
    65 = 1
    66 = 1
    67 = 1
    68 = 2
    69 = 1
    70 = 2
    71 = 2
    72 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
